$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T9").Value = "log10(pCO2) = "
$ws.Range("U9").Value = "(1.482088436399 +  1.302620509340*Z7 +  0.038254786156*T7 +  -0.139692743505*Z7^2 +  -0.037555238467*Z7*T7 +  -0.000769389499*T7^2 +  0.019624637219*Z7^2*T7 +  0.001235261007*Z7*T7^2 +  -0.000657950594*Z7^2*T7^2)"

$ws.Range("AD7").Formula = "=10^(1.482088436399 + 1.302620509340*Z7 + 0.038254786156*T7 + -0.139692743505*Z7^2 + -0.037555238467*Z7*T7 + -0.000769389499*T7^2 + 0.019624637219*Z7^2*T7 + 0.001235261007*Z7*T7^2 + -0.000657950594*Z7^2*T7^2)"

$ws.Columns.Item(20).ColumnWidth = 15.71
$ws.Columns.Item(26).ColumnWidth = 15.71
$ws.Columns.Item(30).ColumnWidth = 15.71
